# "remove in and out"
#
# The "Instruction Set" sheet documents the IN r / OUT r instructions in
# rows 237-248 (A237:L248). Removing that whole block of rows is the
# substance of the commit - everything else (shrunk dimension, shifted
# shared-strings table, the Opcodes sheet's lookup formulas resolving to
# "" for opcodes 0x40-0x45 / 0x48-0x4D, the shrunk conditional-formatting
# range, etc.) is a mechanical consequence of that single structural edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instruction Set")

# Delete rows 237 through 248 inclusive (12 rows: the "IN r" block at
# 237-242 and the "OUT r" block at 243-248). This shifts every following
# row up by 12, recalculates all dependent formulas (e.g. the Opcodes
# sheet's INDEX/MATCH lookups against 'Instruction Set'!I:I, which now
# come back blank for the removed IN/OUT opcodes), and shrinks the
# sheet's used range from A1:L248 down to A1:L236.
$ws.Range("A237:L248").EntireRow.Delete()
